$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: BTC -> BTC
$ws.Range("D2").Value2 = 30557
$ws.Range("E2").Value2 = 593304233415
$ws.Range("F2").Value2 = 7198239265
$ws.Range("G2").Value2 = -0.10056

# Row 3: ETH -> ETH
$ws.Range("D3").Value2 = 1919.94
$ws.Range("E3").Value2 = 230729395706
$ws.Range("F3").Value2 = 4514773379
$ws.Range("G3").Value2 = -0.32974

# Row 4: USDT -> USDT
$ws.Range("D4").Value2 = 0.999861
$ws.Range("E4").Value2 = 83301420286
$ws.Range("F4").Value2 = 11563905746
$ws.Range("G4").Value2 = -0.01214

# Row 5: BNB -> BNB
$ws.Range("D5").Value2 = 245.63
$ws.Range("E5").Value2 = 38277146235
$ws.Range("F5").Value2 = 788571915
$ws.Range("G5").Value2 = 0.81976

# Row 6: USDC -> USDC
$ws.Range("E6").Value2 = 27371288118
$ws.Range("F6").Value2 = 2063985824
$ws.Range("G6").Value2 = -0.06133

# Row 7: XRP -> XRP
$ws.Range("D7").Value2 = 0.478886
$ws.Range("E7").Value2 = 25039187413
$ws.Range("F7").Value2 = 671369535
$ws.Range("G7").Value2 = 1.47385

# Row 8: STETH -> STETH
$ws.Range("D8").Value2 = 1919.32
$ws.Range("E8").Value2 = 14446035044
$ws.Range("F8").Value2 = 1343687
$ws.Range("G8").Value2 = -0.39787

# Row 9: ADA -> ADA
$ws.Range("D9").Value2 = 0.29049
$ws.Range("E9").Value2 = 10178271841
$ws.Range("F9").Value2 = 167533435
$ws.Range("G9").Value2 = 0.66413

# Row 10: DOGE -> DOGE
$ws.Range("D10").Value2 = 0.06725100000000001
$ws.Range("E10").Value2 = 9429648540
$ws.Range("F10").Value2 = 335909511
$ws.Range("G10").Value2 = -1.21056

# Row 11: SOL -> LTC
$ws.Range("B11").Value2 = "LTC"
$ws.Range("C11").Value2 = "Litecoin"
$ws.Range("D11").Value2 = 110.92
$ws.Range("E11").Value2 = 8137896891
$ws.Range("F11").Value2 = 1589312202
$ws.Range("G11").Value2 = 3.86075

# Row 12: TRX -> SOL
$ws.Range("B12").Value2 = "SOL"
$ws.Range("C12").Value2 = "Solana"
$ws.Range("D12").Value2 = 19.04
$ws.Range("E12").Value2 = 7642140711
$ws.Range("F12").Value2 = 354930866
$ws.Range("G12").Value2 = 3.31162

# Row 13: LTC -> TRX
$ws.Range("B13").Value2 = "TRX"
$ws.Range("C13").Value2 = "TRON"
$ws.Range("D13").Value2 = 0.075729
$ws.Range("E13").Value2 = 6805182509
$ws.Range("F13").Value2 = 239364760
$ws.Range("G13").Value2 = -2.49318

# Row 14: DOT -> DOT
$ws.Range("D14").Value2 = 5.29
$ws.Range("E14").Value2 = 6608382738
$ws.Range("F14").Value2 = 111315656
$ws.Range("G14").Value2 = -0.79614

# Row 15: MATIC -> MATIC
$ws.Range("D15").Value2 = 0.671883
$ws.Range("E15").Value2 = 6256783103
$ws.Range("F15").Value2 = 164193079
$ws.Range("G15").Value2 = 1.35651

# Row 16: WBTC -> BCH
$ws.Range("B16").Value2 = "BCH"
$ws.Range("C16").Value2 = "Bitcoin Cash"
$ws.Range("D16").Value2 = 302.06
$ws.Range("E16").Value2 = 5833727668
$ws.Range("F16").Value2 = 1385395110
$ws.Range("G16").Value2 = 2.57902

# Row 17: SHIB -> WBTC
$ws.Range("B17").Value2 = "WBTC"
$ws.Range("C17").Value2 = "Wrapped Bitcoin"
$ws.Range("D17").Value2 = 30579
$ws.Range("E17").Value2 = 4796973137
$ws.Range("F17").Value2 = 49181943
$ws.Range("G17").Value2 = -0.03221

# Row 18: AVAX -> AVAX
$ws.Range("D18").Value2 = 12.99
$ws.Range("E18").Value2 = 4488750007
$ws.Range("F18").Value2 = 104811171
$ws.Range("G18").Value2 = 0.04742

# Row 19: DAI -> SHIB
$ws.Range("B19").Value2 = "SHIB"
$ws.Range("C19").Value2 = "Shiba Inu"
$ws.Range("D19").Value2 = 0.00000757
$ws.Range("E19").Value2 = 4467585055
$ws.Range("F19").Value2 = 77714208
$ws.Range("G19").Value2 = -0.53167

# Row 20: BUSD -> DAI
$ws.Range("B20").Value2 = "DAI"
$ws.Range("C20").Value2 = "Dai"
$ws.Range("D20").Value2 = 0.9994769999999999
$ws.Range("E20").Value2 = 4344063201
$ws.Range("F20").Value2 = 85636279
$ws.Range("G20").Value2 = -0.04935

# Row 21: UNI -> UNI
$ws.Range("D21").Value2 = 5.66
$ws.Range("E21").Value2 = 4259258962
$ws.Range("F21").Value2 = 90815224
$ws.Range("G21").Value2 = 5.29538

# Row 22: BCH -> BUSD
$ws.Range("B22").Value2 = "BUSD"
$ws.Range("C22").Value2 = "Binance USD"
$ws.Range("D22").Value2 = 1
$ws.Range("E22").Value2 = 4133551057
$ws.Range("F22").Value2 = 1358327819
$ws.Range("G22").Value2 = 0.01817

# Row 23: LEO -> LEO
$ws.Range("D23").Value2 = 3.96
$ws.Range("E23").Value2 = 3671042216
$ws.Range("F23").Value2 = 1022295
$ws.Range("G23").Value2 = -2.29074

# Row 24: LINK -> LINK
$ws.Range("D24").Value2 = 6.5
$ws.Range("E24").Value2 = 3375016207
$ws.Range("F24").Value2 = 197370895
$ws.Range("G24").Value2 = 3.93766

# Row 25: TUSD -> TUSD
$ws.Range("D25").Value2 = 0.999413
$ws.Range("E25").Value2 = 3054753793
$ws.Range("F25").Value2 = 776243138
$ws.Range("G25").Value2 = -0.21021

# Row 26: XMR -> XMR
$ws.Range("D26").Value2 = 164.92
$ws.Range("E26").Value2 = 2989837150
$ws.Range("F26").Value2 = 49900901
$ws.Range("G26").Value2 = -2.25734

# Row 27: ATOM -> XLM
$ws.Range("B27").Value2 = "XLM"
$ws.Range("C27").Value2 = "Stellar"
$ws.Range("D27").Value2 = 0.107042
$ws.Range("E27").Value2 = 2903629265
$ws.Range("F27").Value2 = 59296217
$ws.Range("G27").Value2 = -0.59162

# Row 28: OKB -> ETC
$ws.Range("B28").Value2 = "ETC"
$ws.Range("C28").Value2 = "Ethereum Classic"
$ws.Range("D28").Value2 = 20.32
$ws.Range("E28").Value2 = 2878173197
$ws.Range("F28").Value2 = 307572866
$ws.Range("G28").Value2 = -5.09377

# Row 29: ETC -> ATOM
$ws.Range("B29").Value2 = "ATOM"
$ws.Range("C29").Value2 = "Cosmos Hub"
$ws.Range("D29").Value2 = 9.48
$ws.Range("E29").Value2 = 2778027710
$ws.Range("F29").Value2 = 78016681
$ws.Range("G29").Value2 = 1.17496

# Row 30: XLM -> OKB
$ws.Range("B30").Value2 = "OKB"
$ws.Range("C30").Value2 = "OKB"
$ws.Range("D30").Value2 = 44.39
$ws.Range("E30").Value2 = 2662361327
$ws.Range("F30").Value2 = 2753550
$ws.Range("G30").Value2 = 0.0523

# Row 31: TON -> TON
$ws.Range("D31").Value2 = 1.4
$ws.Range("E31").Value2 = 2060300592
$ws.Range("F31").Value2 = 8861564
$ws.Range("G31").Value2 = 2.11258

# Row 32: ICP -> LDO
$ws.Range("B32").Value2 = "LDO"
$ws.Range("C32").Value2 = "Lido DAO"
$ws.Range("D32").Value2 = 2.11
$ws.Range("E32").Value2 = 1859321263
$ws.Range("F32").Value2 = 33726878
$ws.Range("G32").Value2 = 0.13053

# Row 33: FIL -> ICP
$ws.Range("B33").Value2 = "ICP"
$ws.Range("C33").Value2 = "Internet Computer"
$ws.Range("D33").Value2 = 4.16
$ws.Range("E33").Value2 = 1819294894
$ws.Range("F33").Value2 = 11382004
$ws.Range("G33").Value2 = -0.50182

# Row 34: LDO -> FIL
$ws.Range("B34").Value2 = "FIL"
$ws.Range("C34").Value2 = "Filecoin"
$ws.Range("D34").Value2 = 4.06
$ws.Range("E34").Value2 = 1756180896
$ws.Range("F34").Value2 = 114486893
$ws.Range("G34").Value2 = 1.41531

# Row 35: HBAR -> HBAR
$ws.Range("D35").Value2 = 0.050221
$ws.Range("E35").Value2 = 1619633027
$ws.Range("F35").Value2 = 14670162
$ws.Range("G35").Value2 = -0.61434

# Row 36: APT -> QNT
$ws.Range("B36").Value2 = "QNT"
$ws.Range("C36").Value2 = "Quant"
$ws.Range("D36").Value2 = 111.03
$ws.Range("E36").Value2 = 1618370663
$ws.Range("F36").Value2 = 13324579
$ws.Range("G36").Value2 = 0.34708

# Row 37: QNT -> APT
$ws.Range("B37").Value2 = "APT"
$ws.Range("C37").Value2 = "Aptos"
$ws.Range("D37").Value2 = 7.28
$ws.Range("E37").Value2 = 1521084062
$ws.Range("F37").Value2 = 50360207
$ws.Range("G37").Value2 = 0.64051

# Row 38: CRO -> CRO
$ws.Range("D38").Value2 = 0.056742
$ws.Range("E38").Value2 = 1483568361
$ws.Range("F38").Value2 = 5532221
$ws.Range("G38").Value2 = -0.17889

# Row 39: ARB -> VET
$ws.Range("B39").Value2 = "VET"
$ws.Range("C39").Value2 = "VeChain"
$ws.Range("D39").Value2 = 0.02028219
$ws.Range("E39").Value2 = 1475395441
$ws.Range("F39").Value2 = 51341100
$ws.Range("G39").Value2 = -3.59616

# Row 40: NEAR -> ARB
$ws.Range("B40").Value2 = "ARB"
$ws.Range("C40").Value2 = "Arbitrum"
$ws.Range("D40").Value2 = 1.14
$ws.Range("E40").Value2 = 1453902685
$ws.Range("F40").Value2 = 130353211
$ws.Range("G40").Value2 = -1.36545

# Row 41: VET -> NEAR
$ws.Range("B41").Value2 = "NEAR"
$ws.Range("C41").Value2 = "NEAR Protocol"
$ws.Range("D41").Value2 = 1.46
$ws.Range("E41").Value2 = 1355731029
$ws.Range("F41").Value2 = 58594434
$ws.Range("G41").Value2 = 5.41721

# Row 42: AAVE -> AAVE
$ws.Range("D42").Value2 = 71.56999999999999
$ws.Range("E42").Value2 = 1038801939
$ws.Range("F42").Value2 = 108421113
$ws.Range("G42").Value2 = 5.23283

# Row 43: GRT -> FRAX
$ws.Range("B43").Value2 = "FRAX"
$ws.Range("C43").Value2 = "Frax"
$ws.Range("D43").Value2 = 0.999134
$ws.Range("E43").Value2 = 1002564181
$ws.Range("F43").Value2 = 3824535
$ws.Range("G43").Value2 = -0.05847

# Row 44: STX -> GRT
$ws.Range("B44").Value2 = "GRT"
$ws.Range("C44").Value2 = "The Graph"
$ws.Range("D44").Value2 = 0.109234
$ws.Range("E44").Value2 = 996932030
$ws.Range("F44").Value2 = 47827544
$ws.Range("G44").Value2 = 3.17619

# Row 45: ALGO -> USDP
$ws.Range("B45").Value2 = "USDP"
$ws.Range("C45").Value2 = "Pax Dollar"
$ws.Range("D45").Value2 = 1
$ws.Range("E45").Value2 = 987674495
$ws.Range("F45").Value2 = 102178234
$ws.Range("G45").Value2 = 0.05945

# Row 46: FRAX -> RETH
$ws.Range("B46").Value2 = "RETH"
$ws.Range("C46").Value2 = "Rocket Pool ETH"
$ws.Range("D46").Value2 = 2065.81
$ws.Range("E46").Value2 = 954468568
$ws.Range("F46").Value2 = 1436531
$ws.Range("G46").Value2 = -0.34163

# Row 47: USDP -> BSV
$ws.Range("B47").Value2 = "BSV"
$ws.Range("C47").Value2 = "Bitcoin SV"
$ws.Range("D47").Value2 = 49.48
$ws.Range("E47").Value2 = 954284320
$ws.Range("F47").Value2 = 94396984
$ws.Range("G47").Value2 = -1.79127

# Row 48: RETH -> STX
$ws.Range("B48").Value2 = "STX"
$ws.Range("C48").Value2 = "Stacks"
$ws.Range("D48").Value2 = 0.683177
$ws.Range("E48").Value2 = 951441071
$ws.Range("F48").Value2 = 18053415
$ws.Range("G48").Value2 = 0.61628

# Row 49: EGLD -> ALGO
$ws.Range("B49").Value2 = "ALGO"
$ws.Range("C49").Value2 = "Algorand"
$ws.Range("D49").Value2 = 0.123727
$ws.Range("E49").Value2 = 915147177
$ws.Range("F49").Value2 = 44120063
$ws.Range("G49").Value2 = 1.33658

# Row 50: FTM -> EGLD
$ws.Range("B50").Value2 = "EGLD"
$ws.Range("C50").Value2 = "MultiversX"
$ws.Range("D50").Value2 = 35.12
$ws.Range("E50").Value2 = 901038904
$ws.Range("F50").Value2 = 8997558
$ws.Range("G50").Value2 = -0.25455

# Row 51: OP -> OP
$ws.Range("D51").Value2 = 1.33
$ws.Range("E51").Value2 = 859639200
$ws.Range("F51").Value2 = 70041934
$ws.Range("G51").Value2 = -1.69323
